$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = "Caught"
$ws.Range("E2").Value = " Dwaine Pretorius"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1
$ws.Range("N2").Value = " Trent Boult"

# --- Row 3 ---
$ws.Range("D3").Value = "Bowled"
$ws.Range("E3").Value = " Tabraiz Shamsi"
$ws.Range("K3").Value = 11
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = "Caught"
$ws.Range("N3").Value = " Adam Milne"

# --- Row 4 ---
$ws.Range("B4").Value = 44
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = "LBW"
$ws.Range("E4").Value = " Kagiso Rabada"
$ws.Range("K4").Value = 17
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = "Caught"
$ws.Range("N4").Value = " Tim Southee"

# --- Row 5 ---
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "LBW"
$ws.Range("E5").Value = " Kagiso Rabada"
$ws.Range("K5").Value = 26
$ws.Range("L5").Value = 9
$ws.Range("N5").Value = " Adam Milne"

# --- Row 6 ---
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("K6").Value = 13
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = "Caught"
$ws.Range("N6").Value = " Trent Boult"

# --- Row 7 ---
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "Caught"
$ws.Range("E7").Value = " Dwaine Pretorius"
$ws.Range("K7").Value = 10
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = "Caught"
$ws.Range("N7").Value = " Mitchell Santner"

# --- Row 8 ---
$ws.Range("B8").Value = 24
$ws.Range("C8").Value = 9
$ws.Range("D8").Value = "Bowled"
$ws.Range("E8").Value = " Dwaine Pretorius"
$ws.Range("K8").Value = 6
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = "Bowled"
$ws.Range("N8").Value = " Mitchell Santner"

# --- Row 9 ---
$ws.Range("B9").Value = 59
$ws.Range("C9").Value = 24
$ws.Range("E9").Value = " Tabraiz Shamsi"
$ws.Range("K9").Value = 44
$ws.Range("L9").Value = 17
$ws.Range("M9").Value = "Bowled"
$ws.Range("N9").Value = " Tim Southee"

# --- Row 10 ---
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = "LBW"
$ws.Range("K10").Value = 47
$ws.Range("L10").Value = 13
$ws.Range("M10").Value = "Caught"
$ws.Range("N10").Value = " Mitchell Santner"

# --- Row 11 ---
$ws.Range("B11").Value = 24
$ws.Range("C11").Value = 9
$ws.Range("D11").Value = "NOT OUT"
$ws.Range("E11").Value = " "
$ws.Range("K11").Value = 24
$ws.Range("L11").Value = 12
$ws.Range("M11").Value = "NOT OUT"

# --- Row 12 ---
$ws.Range("B12").Value = 7
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = "LBW"
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = "* NOT OUT"

# --- Row 16 ---
$ws.Range("A16").Value = 200
$ws.Range("C16").Value = "'13.0"
$ws.Range("D16").Value = 78
$ws.Range("J16").Value = 201
$ws.Range("K16").Value = 9
$ws.Range("L16").Value = "'12.5"
$ws.Range("M16").Value = 77

# --- Row 21 ---
$ws.Range("B21").Value = "'2.0"
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 14
$ws.Range("K21").Value = "'2.0"
$ws.Range("L21").Value = 30
$ws.Range("M21").Value = 2
$ws.Range("N21").Value = 15

# --- Row 22 ---
$ws.Range("B22").Value = "'2.0"
$ws.Range("C22").Value = 34
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 17
$ws.Range("K22").Value = "'2.0"
$ws.Range("L22").Value = 30
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 15

# --- Row 23 ---
$ws.Range("B23").Value = "'3.0"
$ws.Range("C23").Value = 38
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 12.67
$ws.Range("K23").Value = "'3.0"
$ws.Range("L23").Value = 38
$ws.Range("M23").Value = 2
$ws.Range("N23").Value = 12.67

# --- Row 24 ---
$ws.Range("B24").Value = "'3.0"
$ws.Range("C24").Value = 58
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 19.33
$ws.Range("K24").Value = "'3.0"
$ws.Range("L24").Value = 42
$ws.Range("M24").Value = 2
$ws.Range("N24").Value = 14

# --- Row 25 ---
$ws.Range("B25").Value = "'3.0"
$ws.Range("C25").Value = 42
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 14
$ws.Range("K25").Value = "'2.5"
$ws.Range("L25").Value = 61
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 24.4
